$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '28.256.57'
$ws.Range("E2").Value = '  +3.54%  '

# Row 3
$ws.Range("D3").Value = '1.782.72'
$ws.Range("E3").Value = '  -0.43%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.005'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.40%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '339.15'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.35%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.001'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.44%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3831'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -2.54%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3432'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.45%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '46.88'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.53%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.149'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -4.00%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07382'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.07%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '23.38'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +7.15%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.001'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.33%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.447'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.33%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.367'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +3.70%  '

# Row 16
$ws.Range("D16").Value = '1.782.59'
$ws.Range("E16").Value = '  -0.25%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001075'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.68%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06677'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.24%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '81.80'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.67%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.001'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.17%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.41'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.04%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.421'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.53%  '

# Row 23
$ws.Range("D23").Value = '28.263.33'
$ws.Range("E23").Value = '  +3.57%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '12.07'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.36%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.366'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.81%  '

# Row 26
$ws.Range("B26").Value = 'EthereumClassic'
$ws.Range("C26").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '20.66'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.41%  '

# Row 27
$ws.Range("B27").Value = 'ImmutableX'
$ws.Range("C27").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.424'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -5.01%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.410'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -4.61%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '154.23'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.59%  '

# Row 30
$ws.Range("D30").Value = '1.983.05'
$ws.Range("E30").Value = '  -0.23%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '134.96'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.34%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.016'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.00%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.087'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.78%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.08890'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.97%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '12.72'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -2.62%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02416'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.43%  '

# Row 37
$ws.Range("E37").Value = '  +0.29%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.358'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.46%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.06387'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.74%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.2163'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.58%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.245'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.13%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.499'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -6.93%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.301'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.05%  '

# Row 44
$ws.Range("B44").Value = 'Frax'
$ws.Range("C44").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.000'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.34%  '

# Row 45
$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '14.06'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.47%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.6282'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.24%  '

# Row 47
$ws.Range("E47").Value = '  +0.36%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '132.87'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.36%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.076'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.79%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.07485'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +5.12%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.205'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +6.55%  '
